$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update msg.Pose.Position.X / .Y values (Goal1 column) per the 2019.2 PreScan update
$ws.Range("B2").Value = 264
$ws.Range("B3").Value = 156

# Move/restore the active selection to B4 (as last saved by the author)
[void]$ws.Range("B4").Select()
